$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D4"  = -8.072999999999999
    "B9"  = 5.478999999999999
    "D9"  = -7.992
    "B18" = 5.137
    "B20" = 6.920999999999999
    "D23" = -7.993
    "D24" = -6.837000000000001
    "D26" = -7.558000000000002
    "B27" = 5.752000000000001
    "D34" = -7.564
    "D35" = -7.877999999999998
    "D48" = -7.892999999999999
    "D52" = -7.979000000000001
    "D66" = -7.486
    "D67" = -7.616
    "B69" = 5.843999999999999
    "B76" = 6.544999999999999
    "D80" = -7.938999999999998
    "B82" = 5.456999999999999
    "D99" = -8.253
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
